$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Cells.Item(2,1).Value = "CFAO MOTORS CI"
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 3
$ws1.Cells.Item(2,4).Value = 2705
$ws1.Cells.Item(2,5).Value = 905
$ws1.Cells.Item(2,6).Value = "🟡 Observer"
$ws1.Cells.Item(2,7).Value = "➖ Neutre"

$ws1.Cells.Item(3,1).Value = "UNIWAX CI"
$ws1.Cells.Item(3,2).Value = 0
$ws1.Cells.Item(3,3).Value = 3
$ws1.Cells.Item(3,4).Value = 2635
$ws1.Cells.Item(3,5).Value = 925
$ws1.Cells.Item(3,6).Value = "🟡 Observer"
$ws1.Cells.Item(3,7).Value = "➖ Neutre"

$ws1.Cells.Item(4,1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(4,2).Value = 0
$ws1.Cells.Item(4,3).Value = 6
$ws1.Cells.Item(4,4).Value = 2484.63
$ws1.Cells.Item(4,5).Value = 106.7
$ws1.Cells.Item(4,6).Value = "🟡 Observer"
$ws1.Cells.Item(4,7).Value = "➖ Neutre"

$ws1.Cells.Item(5,1).Value = "SETAO CI"
$ws1.Cells.Item(5,2).Value = 0
$ws1.Cells.Item(5,3).Value = 3
$ws1.Cells.Item(5,4).Value = 2125
$ws1.Cells.Item(5,5).Value = 710
$ws1.Cells.Item(5,6).Value = "🟡 Observer"
$ws1.Cells.Item(5,7).Value = "➖ Neutre"

$ws1.Cells.Item(6,1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(6,2).Value = 0
$ws1.Cells.Item(6,3).Value = 3
$ws1.Cells.Item(6,4).Value = 2005
$ws1.Cells.Item(6,5).Value = 660
$ws1.Cells.Item(6,6).Value = "🟡 Observer"
$ws1.Cells.Item(6,7).Value = "➖ Neutre"

$ws1.Cells.Item(7,1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(7,2).Value = 0
$ws1.Cells.Item(7,3).Value = 3
$ws1.Cells.Item(7,4).Value = 1854.51
$ws1.Cells.Item(7,5).Value = 618.91
$ws1.Cells.Item(7,6).Value = "🟡 Observer"
$ws1.Cells.Item(7,7).Value = "➖ Neutre"

$ws1.Cells.Item(8,1).Value = "AIR LIQUIDE CI"
$ws1.Cells.Item(8,2).Value = 0
$ws1.Cells.Item(8,3).Value = 3
$ws1.Cells.Item(8,4).Value = 1815
$ws1.Cells.Item(8,5).Value = 625
$ws1.Cells.Item(8,6).Value = "🟡 Observer"
$ws1.Cells.Item(8,7).Value = "➖ Neutre"

$ws1.Cells.Item(9,1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(9,2).Value = 0
$ws1.Cells.Item(9,3).Value = 3
$ws1.Cells.Item(9,4).Value = 1184.71
$ws1.Cells.Item(9,5).Value = 393.46
$ws1.Cells.Item(9,6).Value = "🟡 Observer"
$ws1.Cells.Item(9,7).Value = "➖ Neutre"

$ws1.Cells.Item(10,1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(10,2).Value = 0
$ws1.Cells.Item(10,3).Value = 3
$ws1.Cells.Item(10,4).Value = 1043.96
$ws1.Cells.Item(10,5).Value = 347.58
$ws1.Cells.Item(10,6).Value = "🟡 Observer"
$ws1.Cells.Item(10,7).Value = "➖ Neutre"

$ws1.Cells.Item(11,1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(11,2).Value = 0
$ws1.Cells.Item(11,3).Value = 3
$ws1.Cells.Item(11,4).Value = 993.0599999999999
$ws1.Cells.Item(11,5).Value = 332.11
$ws1.Cells.Item(11,6).Value = "🟡 Observer"
$ws1.Cells.Item(11,7).Value = "➖ Neutre"

$ws1.Cells.Item(12,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(12,2).Value = 0
$ws1.Cells.Item(12,3).Value = 3
$ws1.Cells.Item(12,4).Value = 416.23
$ws1.Cells.Item(12,5).Value = 139.39
$ws1.Cells.Item(12,6).Value = "🟡 Observer"
$ws1.Cells.Item(12,7).Value = "➖ Neutre"

$ws1.Cells.Item(13,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(13,2).Value = 0
$ws1.Cells.Item(13,3).Value = 3
$ws1.Cells.Item(13,4).Value = 402.03
$ws1.Cells.Item(13,5).Value = 134.12
$ws1.Cells.Item(13,6).Value = "🟡 Observer"
$ws1.Cells.Item(13,7).Value = "➖ Neutre"

$ws1.Cells.Item(14,1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(14,2).Value = 0
$ws1.Cells.Item(14,3).Value = 3
$ws1.Cells.Item(14,4).Value = 389.12
$ws1.Cells.Item(14,5).Value = 131
$ws1.Cells.Item(14,6).Value = "🟡 Observer"
$ws1.Cells.Item(14,7).Value = "➖ Neutre"

$ws1.Cells.Item(15,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(15,2).Value = 0
$ws1.Cells.Item(15,3).Value = 3
$ws1.Cells.Item(15,4).Value = 382.43
$ws1.Cells.Item(15,5).Value = 128.75
$ws1.Cells.Item(15,6).Value = "🟡 Observer"
$ws1.Cells.Item(15,7).Value = "➖ Neutre"

$ws1.Cells.Item(16,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(16,2).Value = 0
$ws1.Cells.Item(16,3).Value = 3
$ws1.Cells.Item(16,4).Value = 369.39
$ws1.Cells.Item(16,5).Value = 123.05
$ws1.Cells.Item(16,6).Value = "🟡 Observer"
$ws1.Cells.Item(16,7).Value = "➖ Neutre"

$ws1.Cells.Item(17,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(17,2).Value = 0
$ws1.Cells.Item(17,3).Value = 3
$ws1.Cells.Item(17,4).Value = 323.72
$ws1.Cells.Item(17,5).Value = 107.78
$ws1.Cells.Item(17,6).Value = "🟡 Observer"
$ws1.Cells.Item(17,7).Value = "➖ Neutre"

$ws1.Cells.Item(18,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(18,2).Value = 0
$ws1.Cells.Item(18,3).Value = 3
$ws1.Cells.Item(18,4).Value = 282.6
$ws1.Cells.Item(18,5).Value = 93.90000000000001
$ws1.Cells.Item(18,6).Value = "🟡 Observer"
$ws1.Cells.Item(18,7).Value = "➖ Neutre"

$ws1.Cells.Item(19,1).Value = "BRVM - INDUSTRIE                (**)"
$ws1.Cells.Item(19,2).Value = 0
$ws1.Cells.Item(19,3).Value = 1
$ws1.Cells.Item(19,4).Value = 244.24
$ws1.Cells.Item(19,5).Value = 244.24
$ws1.Cells.Item(19,6).Value = "🟡 Observer"
$ws1.Cells.Item(19,7).Value = "➖ Neutre"

$ws1.Cells.Item(20,1).Value = "BRVM - CONSOMMATION DE BASE             (**)"
$ws1.Cells.Item(20,2).Value = 0
$ws1.Cells.Item(20,3).Value = 1
$ws1.Cells.Item(20,4).Value = 204.59
$ws1.Cells.Item(20,5).Value = 204.59
$ws1.Cells.Item(20,6).Value = "🟡 Observer"
$ws1.Cells.Item(20,7).Value = "➖ Neutre"

$ws1.Cells.Item(21,1).Value = "BRVM-PRINCIPAL                  (**)"
$ws1.Cells.Item(21,2).Value = 0
$ws1.Cells.Item(21,3).Value = 1
$ws1.Cells.Item(21,4).Value = 192.37
$ws1.Cells.Item(21,5).Value = 192.37
$ws1.Cells.Item(21,6).Value = "🟡 Observer"
$ws1.Cells.Item(21,7).Value = "➖ Neutre"

$ws1.Cells.Item(22,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(22,2).Value = 2
$ws1.Cells.Item(22,3).Value = 0
$ws1.Cells.Item(22,4).Value = 14.39
$ws1.Cells.Item(22,5).Value = 6.94
$ws1.Cells.Item(22,6).Value = "🟡 Observer"
$ws1.Cells.Item(22,7).Value = "👀 À surveiller"

$ws1.Cells.Item(23,1).Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Cells.Item(23,2).Value = 2
$ws1.Cells.Item(23,3).Value = 0
$ws1.Cells.Item(23,4).Value = 11.37
$ws1.Cells.Item(23,5).Value = 7.2
$ws1.Cells.Item(23,6).Value = "🟡 Observer"
$ws1.Cells.Item(23,7).Value = "➖ Neutre"

$ws1.Cells.Item(24,1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(24,2).Value = 2
$ws1.Cells.Item(24,3).Value = 1
$ws1.Cells.Item(24,4).Value = 7.21
$ws1.Cells.Item(24,5).Value = 7.2
$ws1.Cells.Item(24,6).Value = "🟡 Observer"
$ws1.Cells.Item(24,7).Value = "👀 À surveiller"

$ws1.Cells.Item(25,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(25,2).Value = 1
$ws1.Cells.Item(25,3).Value = 0
$ws1.Cells.Item(25,4).Value = 7.04
$ws1.Cells.Item(25,5).Value = 7.04
$ws1.Cells.Item(25,6).Value = "🟡 Observer"
$ws1.Cells.Item(25,7).Value = "➖ Neutre"

$ws1.Cells.Item(26,1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(26,2).Value = 2
$ws1.Cells.Item(26,3).Value = 0
$ws1.Cells.Item(26,4).Value = 6.12
$ws1.Cells.Item(26,5).Value = 3.27
$ws1.Cells.Item(26,6).Value = "🟢 Achat"
$ws1.Cells.Item(26,7).Value = "✅ Renforcer"

$ws1.Cells.Item(27,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(27,2).Value = 1
$ws1.Cells.Item(27,3).Value = 0
$ws1.Cells.Item(27,4).Value = 5.88
$ws1.Cells.Item(27,5).Value = 5.88
$ws1.Cells.Item(27,6).Value = "🟡 Observer"
$ws1.Cells.Item(27,7).Value = "➖ Neutre"

$ws1.Cells.Item(28,1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Cells.Item(28,2).Value = 1
$ws1.Cells.Item(28,3).Value = 0
$ws1.Cells.Item(28,4).Value = 4.95
$ws1.Cells.Item(28,5).Value = 4.95
$ws1.Cells.Item(28,6).Value = "🟡 Observer"
$ws1.Cells.Item(28,7).Value = "➖ Neutre"

$ws1.Cells.Item(29,1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(29,2).Value = 1
$ws1.Cells.Item(29,3).Value = 0
$ws1.Cells.Item(29,4).Value = 4.53
$ws1.Cells.Item(29,5).Value = 4.53
$ws1.Cells.Item(29,6).Value = "🟡 Observer"
$ws1.Cells.Item(29,7).Value = "➖ Neutre"

$ws1.Cells.Item(30,1).Value = "BANK OF AFRICA CI (BOAC)"
$ws1.Cells.Item(30,2).Value = 1
$ws1.Cells.Item(30,3).Value = 0
$ws1.Cells.Item(30,4).Value = 4.27
$ws1.Cells.Item(30,5).Value = 4.27
$ws1.Cells.Item(30,6).Value = "🟡 Observer"
$ws1.Cells.Item(30,7).Value = "➖ Neutre"

$ws1.Cells.Item(31,1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Cells.Item(31,2).Value = 1
$ws1.Cells.Item(31,3).Value = 0
$ws1.Cells.Item(31,4).Value = 3.98
$ws1.Cells.Item(31,5).Value = 3.98
$ws1.Cells.Item(31,6).Value = "🟡 Observer"
$ws1.Cells.Item(31,7).Value = "➖ Neutre"

$ws1.Cells.Item(32,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(32,2).Value = 1
$ws1.Cells.Item(32,3).Value = 0
$ws1.Cells.Item(32,4).Value = 3.01
$ws1.Cells.Item(32,5).Value = 3.01
$ws1.Cells.Item(32,6).Value = "🟡 Observer"
$ws1.Cells.Item(32,7).Value = "➖ Neutre"

$ws1.Cells.Item(33,1).Value = "TOTAL"
$ws1.Cells.Item(33,2).Value = 0
$ws1.Cells.Item(33,3).Value = 3
$ws1.Cells.Item(33,4).Value = 0
$ws1.Cells.Item(33,5).Value = 0
$ws1.Cells.Item(33,6).Value = "🟡 Observer"
$ws1.Cells.Item(33,7).Value = "➖ Neutre"

$ws1.Cells.Item(34,1).Value = "NESTLE CI (NTLC)"
$ws1.Cells.Item(34,2).Value = 0
$ws1.Cells.Item(34,3).Value = 1
$ws1.Cells.Item(34,4).Value = -1.41
$ws1.Cells.Item(34,5).Value = -1.41
$ws1.Cells.Item(34,6).Value = "🟡 Observer"
$ws1.Cells.Item(34,7).Value = "➖ Neutre"

$ws1.Cells.Item(35,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(35,2).Value = 0
$ws1.Cells.Item(35,3).Value = 1
$ws1.Cells.Item(35,4).Value = -1.54
$ws1.Cells.Item(35,5).Value = -1.54
$ws1.Cells.Item(35,6).Value = "🟡 Observer"
$ws1.Cells.Item(35,7).Value = "➖ Neutre"

$ws1.Cells.Item(36,1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(36,2).Value = 0
$ws1.Cells.Item(36,3).Value = 1
$ws1.Cells.Item(36,4).Value = -2.22
$ws1.Cells.Item(36,5).Value = -2.22
$ws1.Cells.Item(36,6).Value = "🟡 Observer"
$ws1.Cells.Item(36,7).Value = "➖ Neutre"

$ws1.Cells.Item(37,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(37,2).Value = 0
$ws1.Cells.Item(37,3).Value = 1
$ws1.Cells.Item(37,4).Value = -2.74
$ws1.Cells.Item(37,5).Value = -2.74
$ws1.Cells.Item(37,6).Value = "🟡 Observer"
$ws1.Cells.Item(37,7).Value = "➖ Neutre"

$ws1.Cells.Item(38,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(38,2).Value = 0
$ws1.Cells.Item(38,3).Value = 2
$ws1.Cells.Item(38,4).Value = -5.86
$ws1.Cells.Item(38,5).Value = -4.17
$ws1.Cells.Item(38,6).Value = "🟡 Observer"
$ws1.Cells.Item(38,7).Value = "➖ Neutre"

$ws1.Cells.Item(39,1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(39,2).Value = 0
$ws1.Cells.Item(39,3).Value = 2
$ws1.Cells.Item(39,4).Value = -7.4
$ws1.Cells.Item(39,5).Value = -4.17
$ws1.Cells.Item(39,6).Value = "🟡 Observer"
$ws1.Cells.Item(39,7).Value = "➖ Neutre"

$ws1.Cells.Item(40,1).Value = "BICI CI (BICC)"
$ws1.Cells.Item(40,2).Value = 0
$ws1.Cells.Item(40,3).Value = 1
$ws1.Cells.Item(40,4).Value = -7.47
$ws1.Cells.Item(40,5).Value = -7.47
$ws1.Cells.Item(40,6).Value = "🟡 Observer"
$ws1.Cells.Item(40,7).Value = "➖ Neutre"

$ws1.Cells.Item(41,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(41,2).Value = 0
$ws1.Cells.Item(41,3).Value = 2
$ws1.Cells.Item(41,4).Value = -14.75
$ws1.Cells.Item(41,5).Value = -7.36
$ws1.Cells.Item(41,6).Value = "🟡 Observer"
$ws1.Cells.Item(41,7).Value = "➖ Neutre"

$ws1.Cells.Item(42,1).Value = "UNILEVER CI (UNLC)"
$ws1.Cells.Item(42,2).Value = 0
$ws1.Cells.Item(42,3).Value = 3
$ws1.Cells.Item(42,4).Value = -22.49
$ws1.Cells.Item(42,5).Value = -7.5
$ws1.Cells.Item(42,6).Value = "🔴 Vente"
$ws1.Cells.Item(42,7).Value = "⚠️ Risque de décrochage"

# Remove now-obsolete trailing rows 43 and 44 (table shrank by two rows)
$ws1.Range("A44:G44").EntireRow.Delete()
$ws1.Range("A43:G43").EntireRow.Delete()

$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Cells.Item(2,2).Value = 484265.68
$ws2.Cells.Item(3,2).Value = 100400
$ws2.Cells.Item(4,2).Value = 93175
$ws2.Cells.Item(5,2).Value = 52716.05
$ws2.Cells.Item(6,2).Value = 45241.6
$ws2.Cells.Item(7,2).Value = 36940.77
$ws2.Cells.Item(8,2).Value = 34917.5
$ws2.Cells.Item(9,2).Value = 12021.55
$ws2.Cells.Item(10,2).Value = 8890.709999999999
$ws2.Cells.Item(11,2).Value = 7907.34
